$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: who filled this in / which week ---
$ws.Range("C1").Value = "Richard Dobson"
$ws.Range("E1").Value = 7

# --- Task rows (Stage / Task / Estimated / Hours Spent) ---
$ws.Range("A3").Value = "Project Build"
$ws.Range("B3").Value = "Finish work on metatdata array populate function"
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 5

$ws.Range("A4").Value = "Project Build"
$ws.Range("B4").Value = "Iteration Review"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1

$ws.Range("A5").Value = "Project Build"
$ws.Range("B5").Value = "Gather any new requirements, commence work on next iteration"
$ws.Range("C5").Value = 14
$ws.Range("D5").Value = 14

# --- Cumulative total label (formula in D14 recalculates automatically) ---
$ws.Range("A14").Value = "Cumulative Total: 140"

# --- Column B a bit wider to fit the longer task text ---
$ws.Columns.Item(2).ColumnWidth = 48.75

# --- Last selected cell when the file was saved ---
$ws.Range("I15").Select()
